$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text in B2 (replaces the old "인재개발부" intro
# with the new "드래곤볼의 손오공" intro)
$ws.Range("B2").Value = "안녕하세요. 드래곤볼의 손오공입니다! "

# Move the active selection to B2 (was B11)
$ws.Range("B2").Select()

# Best-effort: remember the workbook window's on-screen position, matching
# the author's recorded xWindow/yWindow for this save.
$excel.Windows.Item(1).Left = 27900
$excel.Windows.Item(1).Top = 4640
